# CCS-7, SP-30: Updated proteomics handlers to expect format template from Rolf
# - "Value Unit" description gains a new allowed unit ("fmol/ug protein digest")
# - Description column (C) on openbis-metadata is widened to fit the longer text
# - Active selection on openbis-metadata moves up one row (C8 -> C7)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("openbis-metadata")

# Extend the allowed "Value Unit" text in the Description column (row 6).
$ws.Range("C6").Value = "One of mM, uM, Percent, RatioT1, RatioCs, or AU, Dimensionless, fmol/ug protein digest"

# Widen column C (Description) so the longer text fits - target stored width is 90.
# COM's ColumnWidth is in characters and gets ~0.714 of padding added when the
# file is written back out, so dial it in to land exactly on 90.
$ws.Columns.Item(3).ColumnWidth = 89.28571428571429

# Move the selection highlight from C8 to C7.
$ws.Range("C7").Select() | Out-Null
